# Auto-generated edit script: updates column F numeric values across 4 sheets
# per the commit diff ("Update gh-pages to output generated at 456a3b4")
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 6507   # was 6502
$ws.Range("F4").Value = 749   # was 748
$ws.Range("F6").Value = 92   # was 89
$ws.Range("F7").Value = 575   # was 571
$ws.Range("F9").Value = 31   # was 29
$ws.Range("F10").Value = 755   # was 749
$ws.Range("F11").Value = 1227   # was 1225
$ws.Range("F12").Value = 16   # was 14
$ws.Range("F13").Value = 92   # was 89
$ws.Range("F14").Value = 207   # was 206
$ws.Range("F15").Value = 472   # was 467
$ws.Range("F21").Value = 414   # was 411
$ws.Range("F24").Value = 180   # was 177
$ws.Range("F25").Value = 2252   # was 2248
$ws.Range("F30").Value = 3669   # was 3665
$ws.Range("F32").Value = 666   # was 662

# Sheet: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 203   # was 202
$ws.Range("F6").Value = 28   # was 27
$ws.Range("F7").Value = 724   # was 721
$ws.Range("F12").Value = 118   # was 117
$ws.Range("F20").Value = 4100   # was 4098
$ws.Range("F25").Value = 206   # was 205
$ws.Range("F29").Value = 217   # was 215
$ws.Range("F33").Value = 1682   # was 1681
$ws.Range("F34").Value = 25   # was 24

# Sheet: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 48   # was 47
$ws.Range("F4").Value = 1212   # was 1210
$ws.Range("F6").Value = 1590   # was 1588
$ws.Range("F10").Value = 858   # was 851

# Sheet: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 48   # was 47
$ws.Range("F4").Value = 1212   # was 1210
$ws.Range("F5").Value = 1590   # was 1588
$ws.Range("F8").Value = 858   # was 851
$ws.Range("F9").Value = 6507   # was 6502
$ws.Range("F10").Value = 28   # was 27
$ws.Range("F11").Value = 749   # was 748
$ws.Range("F12").Value = 724   # was 721
$ws.Range("F13").Value = 92   # was 89
$ws.Range("F14").Value = 575   # was 571
$ws.Range("F16").Value = 31   # was 29
$ws.Range("F17").Value = 755   # was 749
$ws.Range("F19").Value = 118   # was 117
$ws.Range("F20").Value = 118   # was 117
$ws.Range("F23").Value = 1227   # was 1225
$ws.Range("F24").Value = 16   # was 14
$ws.Range("F25").Value = 92   # was 89
$ws.Range("F26").Value = 207   # was 206
$ws.Range("F35").Value = 414   # was 411
$ws.Range("F39").Value = 180   # was 177
$ws.Range("F40").Value = 2252   # was 2248
$ws.Range("F41").Value = 217   # was 215
$ws.Range("F43").Value = 1682   # was 1681
$ws.Range("F44").Value = 1682   # was 1681
$ws.Range("F47").Value = 3669   # was 3665
$ws.Range("F51").Value = 666   # was 662
